$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  @(0.00158971523361902,0.0000691180536356096,0.0,0.0,0.0,0.999792645839093,0.442977605750622,0.999308819463644,0.000207354160906829,0.0,0.996474979264584,0.0000691180536356096,0.000552944429084877,0.000276472214542438,0.000207354160906829,0.000483826375449267,0.000207354160906829,0.000898534697262925,0.000345590268178048,0.000276472214542438,0.996405861210948,0.915883328725463,0.000967652750898535),
  @(0.997649986176389,0.999723527785458,0.999239701410008,0.000276472214542438,0.999792645839093,0.0,0.0000691180536356096,0.000276472214542438,0.000483826375449267,0.99149847940282,0.00138236107271219,0.998755875034559,0.000276472214542438,0.998963229195466,0.000760298589991706,0.000829416643627315,0.000552944429084877,0.0,0.000138236107271219,0.999032347249102,0.000483826375449267,0.000138236107271219,0.998686756980923),
  @(0.0000691180536356096,0.0,0.0,0.0,0.0000691180536356096,0.0000691180536356096,0.548106165330384,0.000414708321813658,0.000345590268178048,0.0,0.000967652750898535,0.000552944429084877,0.0,0.0000691180536356096,0.0,0.0000691180536356096,0.000483826375449267,0.999032347249102,0.999447055570915,0.000345590268178048,0.00304119435996682,0.080729886646392,0.000138236107271219),
  @(0.000414708321813658,0.000138236107271219,0.000552944429084877,0.999516173624551,0.000138236107271219,0.0000691180536356096,0.000138236107271219,0.0,0.998755875034559,0.00836328448990876,0.000829416643627315,0.000207354160906829,0.999101465302737,0.000207354160906829,0.998963229195466,0.998548520873652,0.998548520873652,0.0000691180536356096,0.0000691180536356096,0.000207354160906829,0.0,0.0000691180536356096,0.0)
)

for ($i = 0; $i -lt 4; $i++) {
  $row = $i + 2
  for ($j = 0; $j -lt 23; $j++) {
    $col = $j + 2
    $ws.Cells.Item($row, $col).Value = $values[$i][$j]
  }
}
